$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force-write a D-column (Price) value as literal text, then restore
# the default "Normal" style so the cell style is unchanged (Excel auto-
# coerces bare numeric-looking strings assigned via .Value into numbers,
# which would lose formatting like trailing zeros or the thousands-dot style
# used in this sheet, e.g. "26.196.45").
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '26.196.45'
$ws.Range("E2").Value = '  -0.30%  '
Set-TextValue $ws.Range("D3") '1.680.50'
$ws.Range("E3").Value = '  +0.10%  '
Set-TextValue $ws.Range("D5") '216.35'
$ws.Range("E5").Value = '  -0.71%  '
Set-TextValue $ws.Range("D6") '0.5261'
$ws.Range("E6").Value = '  -1.52%  '
$ws.Range("E7").Value = '  -0.03%  '
Set-TextValue $ws.Range("D8") '0.2690'
$ws.Range("E8").Value = '  +0.22%  '
$ws.Range("E9").Value = '  -1.55%  '
Set-TextValue $ws.Range("D10") '21.42'
$ws.Range("E10").Value = '  -2.27%  '
Set-TextValue $ws.Range("D11") '0.07619'
Set-TextValue $ws.Range("D12") '1.687.90'
$ws.Range("E12").Value = '  +0.47%  '
Set-TextValue $ws.Range("D13") '4.520'
$ws.Range("E13").Value = '  +0.00%  '
Set-TextValue $ws.Range("D14") '0.5735'
$ws.Range("E14").Value = '  -0.69%  '
Set-TextValue $ws.Range("D15") '0.000008229'
$ws.Range("E15").Value = '  -2.90%  '
Set-TextValue $ws.Range("D16") '66.09'
$ws.Range("E16").Value = '  +2.10%  '
Set-TextValue $ws.Range("D17") '26.220.89'
$ws.Range("E17").Value = '  -0.35%  '
$ws.Range("E18").Value = '  -0.05%  '
Set-TextValue $ws.Range("D19") '4.867'
$ws.Range("E19").Value = '  -0.68%  '
Set-TextValue $ws.Range("D20") '10.74'
$ws.Range("E20").Value = '  -1.13%  '
Set-TextValue $ws.Range("D21") '189.64'
$ws.Range("E21").Value = '  -0.40%  '
Set-TextValue $ws.Range("D22") '6.233'
$ws.Range("E22").Value = '  +0.48%  '
Set-TextValue $ws.Range("D23") '1.006'
$ws.Range("E23").Value = '  -0.10%  '
Set-TextValue $ws.Range("D24") '148.83'
$ws.Range("E24").Value = '  +2.13%  '
Set-TextValue $ws.Range("D25") '0.1260'
$ws.Range("E25").Value = '  -0.93%  '
Set-TextValue $ws.Range("D26") '7.728'
$ws.Range("E26").Value = '  -1.20%  '
Set-TextValue $ws.Range("D27") '15.79'
$ws.Range("E27").Value = '  +0.30%  '
Set-TextValue $ws.Range("D28") '0.06341'
$ws.Range("E28").Value = '  -2.04%  '
$ws.Range("E29").Value = '  -0.05%  '
Set-TextValue $ws.Range("D30") '1.316'
$ws.Range("E30").Value = '  -0.23%  '
Set-TextValue $ws.Range("D31") '3.565'
$ws.Range("E31").Value = '  -0.44%  '
Set-TextValue $ws.Range("D32") '3.563'
$ws.Range("E32").Value = '  -0.52%  '
Set-TextValue $ws.Range("D33") '1.679'
$ws.Range("E33").Value = '  +0.99%  '
Set-TextValue $ws.Range("D34") '1.021'
$ws.Range("E34").Value = '  -1.07%  '
Set-TextValue $ws.Range("D35") '0.6102'
$ws.Range("E35").Value = '  -1.12%  '
Set-TextValue $ws.Range("D36") '2.421'
$ws.Range("E36").Value = '  +0.76%  '
$ws.Range("E37").Value = '  +1.51%  '
$ws.Range("E38").Value = '  -1.21%  '
Set-TextValue $ws.Range("D39") '0.01614'
$ws.Range("E39").Value = '  -0.50%  '
Set-TextValue $ws.Range("D40") '1.095.76'
$ws.Range("E40").Value = '  -1.41%  '
Set-TextValue $ws.Range("D41") '0.8832'
$ws.Range("E41").Value = '  +1.50%  '
Set-TextValue $ws.Range("D43") '100.43'
$ws.Range("E43").Value = '  +0.08%  '
Set-TextValue $ws.Range("D44") '1.831.80'
$ws.Range("E44").Value = '  +0.21%  '
Set-TextValue $ws.Range("D45") '0.00000000109'
$ws.Range("E45").Value = '  +0.56%  '
Set-TextValue $ws.Range("D46") '57.36'
$ws.Range("E46").Value = '  +0.37%  '
$ws.Range("E47").Value = '  +0.08%  '
Set-TextValue $ws.Range("D48") '8.070'
$ws.Range("E48").Value = '  -1.33%  '
$ws.Range("E49").Value = '  +0.15%  '
Set-TextValue $ws.Range("D50") '0.4279'
Set-TextValue $ws.Range("D51") '5.993'
$ws.Range("E51").Value = '  -1.34%  '
